$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first data row (ID "H 72") - everything below shifts up one row,
# which also drops the former last row (old row 63) so the sheet ends at row 62.
$ws.Rows.Item(2).Delete()

# Rewrite the data block (rows 2-62, columns A-F) with the corrected/completed
# values (imputed-data-project ground truth), since the row shift alone does
# not reproduce the target values for every cell.
$data = @(
    @("H 1968 Sk 4", -20, 10, -14.8, -6.9, 0.7098100000000001),
    @("H 12640 K XIII/3", -19.8, 11.7, -13.5, -7.5, 0.70917),
    @("RM 2", $null, 14.9, -13.5, -7.2, 0.70909),
    @("RM 8", -19.7, 11.2, -14.2, $null, 0.71055),
    @("RM 9", -18.7, 11, -15.4, -6.4, 0.71067),
    @("RM 14", $null, 12.3, -14.4, -5, 0.71266),
    @("RM 21", -19.8, 15.1, -14.2, -5.7, 0.7099299999999999),
    @("RM 28", -19.9, 12, -15.4, -6.2, 0.71194),
    @("RM 32", -19.5, 15, -13.8, -7.1, 0.7105),
    @("RM 38", -19.9, 15.5, -13.9, -6.6, 0.71122),
    @("RM 42", -20.6, 10.5, -14.5, -6.8, 0.71118),
    @("RM 52 a", -19.8, 11.5, -14.7, -6.1, 0.7110300000000001),
    @("RM 58", -20.8, 11.4, -15.5, -7.9, 0.70989),
    @("RM 78", -18.9, 12, -15.4, -12, 0.70964),
    @("RM 81", $null, 12.5, -14.1, -5.3, 0.70956),
    @("RM 88", -19.9, 12.5, -13.9, -5.3, 0.71011),
    @("RM 90", -19.6, 14.4, -13.1, -5.4, 0.71073),
    @("RM 95", -19.1, 12.5, -15.2, -8.4, 0.7107599999999999),
    @("RM 103", -19.5, 13.5, -15.3, -6.9, 0.7106),
    @("RM 116", -19.4, 11.2, -14.7, -7.3, 0.7098100000000001),
    @("RM 120", -19.6, 11.5, -15.2, -8.5, 0.70941),
    @("RM 125", -20.6, 13.2, -15.5, -6.5, 0.70931),
    @("RM 134", -19, 12.5, -14, -7.2, 0.7101499999999999),
    @("RM 135", -18.9, 12.7, -14.3, -8.699999999999999, 0.7108),
    @("RM 137", -19.5, 12, -14.9, -8.5, 0.70925),
    @("RM 138", -19.3, 13.5, -15.4, -6.1, 0.7092000000000001),
    @("RM 140", $null, 12.2, -13.9, -7, 0.70963),
    @("RM 142a", $null, 12.7, -13.9, -8.1, 0.7094200000000001),
    @("RM 142b", -19.5, 12.2, -13.8, -7.4, 0.71055),
    @("RM 145", -19.5, 10.7, -15.5, -7.1, 0.7116),
    @("RM 146", $null, 12.6, -15.2, -7.8, 0.71028),
    @("RM 156", -19.4, 12, -15.1, -9.5, 0.72961),
    @("RM 158", -19.6, 10.3, -14, -7.6, 0.70933),
    @("RM 159", -19.2, 11.3, -14.4, -7.4, 0.70933),
    @("RM 165", -19.1, 14.3, -14.6, -8.6, 0.71087),
    @("RM 167", -19.8, 12.1, -14.3, -7.1, 0.71078),
    @("RM 170", -19.2, 11.5, -14.7, -9.199999999999999, 0.70956),
    @("RM 173", -19.8, 11.7, -14.7, $null, 0.70929),
    @("RM 178", $null, 11.9, -14.4, -7.9, 0.71102),
    @("RM 186", -18.1, 13.9, -15.9, -10.2, 0.7111499999999999),
    @("RM 193", -19, 12, -14.6, -6.8, 0.7111499999999999),
    @("RM 197", -18.9, 13, -14.1, -8.800000000000001, 0.71152),
    @("RM 202", -18.8, 12.6, -15.3, -8.9, 0.7385699999999999),
    @("RM 207", -19.7, 11.7, -14.9, -7.4, 0.70941),
    @("RM 208", -19.5, 11.7, -14.9, -7.3, 0.71128),
    @("RM 215", -19.8, 12.9, -14.6, $null, 0.71117),
    @("RM 219", -19.3, 12.7, -12.8, -8.300000000000001, 0.71045),
    @("RM 221", -19.9, 11.5, -13.7, -7.2, 0.70924),
    @("RM 232", -19.7, 10.7, -15.6, -8.800000000000001, 0.71064),
    @("RM 233", -20.5, 11.6, -14.5, -7.7, 0.7107599999999999),
    @("SC 5", $null, 10.8, -13.8, -5, 0.70948),
    @("SC 66", -20.3, 10.5, -12.9, -5.7, 0.71211),
    @("SC 92", -17.2, 14.3, -14, -6.3, 0.71157),
    @("SC 101", -20.4, 10, -14.6, $null, 0.71492),
    @("SC 103", -19.2, 11.9, -14.7, -5.7, 0.7108),
    @("SC 105", $null, 11.1, -13.7, -5.9, 0.71146),
    @("SC 119", -19.5, 11.2, -13, -6.8, 0.70951),
    @("SC 120", -19.7, 11.4, -13.6, -5.7, 0.71096),
    @("SC 132", -18.8, 15.3, -13.7, -8.1, 0.70948),
    @("SC 193", -19.9, 10.5, -14.7, $null, 0.71183),
    @("SC 232", -19.5, 10.4, -14.1, $null, 0.7115899999999999)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]

    if ($null -eq $row[1]) { $ws.Cells.Item($r, 2).Value = "" } else { $ws.Cells.Item($r, 2).Value = $row[1] }
    if ($null -eq $row[2]) { $ws.Cells.Item($r, 3).Value = "" } else { $ws.Cells.Item($r, 3).Value = $row[2] }
    if ($null -eq $row[3]) { $ws.Cells.Item($r, 4).Value = "" } else { $ws.Cells.Item($r, 4).Value = $row[3] }
    if ($null -eq $row[4]) { $ws.Cells.Item($r, 5).Value = "" } else { $ws.Cells.Item($r, 5).Value = $row[4] }
    if ($null -eq $row[5]) { $ws.Cells.Item($r, 6).Value = "" } else { $ws.Cells.Item($r, 6).Value = $row[5] }
}
